$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12925.75
$ws.Range("I62").Value = 13700
$ws.Range("J62").Value = 12667.667
$ws.Range("K62").Value = 13700
$ws.Range("L62").Value = 12667.667
$ws.Range("M62").Value = -13076
$ws.Range("N62").Value = -13915.667

$ws.Range("H65").Value = 12925.75
$ws.Range("I65").Value = 13700
$ws.Range("J65").Value = 12667.667
$ws.Range("K65").Value = 68500
$ws.Range("L65").Value = 63338.335
$ws.Range("M65").Value = -65380
$ws.Range("N65").Value = -69578.33499999999

$ws.Range("H112").Value = 1186.6666
$ws.Range("J112").Value = 1206.1794
$ws.Range("L112").Value = 3618.5382
$ws.Range("N112").Value = -5834.5382

$ws.Range("H132").Value = 2612.6438
$ws.Range("I132").Value = 1795.8636
$ws.Range("K132").Value = 5387.5908
$ws.Range("M132").Value = -2857.5908

$ws.Range("H137").Value = 76925690
$ws.Range("I137").Value = 1000000000
$ws.Range("J137").Value = 2829.0833
$ws.Range("K137").Value = 3000000000
$ws.Range("L137").Value = 8487.249899999999
$ws.Range("M137").Value = -2999997450
$ws.Range("N137").Value = -13587.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 9236.799999999999
$ws.Range("I31").Value = 8125
$ws.Range("J31").Value = 13684
$ws.Range("K31").Value = 8125
$ws.Range("L31").Value = 13684
$ws.Range("M31").Value = -7831
$ws.Range("N31").Value = -14272

$ws.Range("H32").Value = 1114.3269
$ws.Range("I32").Value = 1114.3269
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1114.3269
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = -827.3269
$ws.Range("M32").ClearContents()

$ws.Range("H45").Value = 2416.5833
$ws.Range("I45").Value = 2369.4783
$ws.Range("K45").Value = 2369.4783
$ws.Range("M45").Value = -1992.4783

$ws.Range("H61").Value = 3130.7778
$ws.Range("I61").Value = 2273.5386
$ws.Range("K61").Value = 2273.5386
$ws.Range("M61").Value = -2061.5386

$ws.Range("H75").Value = 76839.164
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 76839.164
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = 76839.164
$ws.Range("N75").Value = -78587.164
$ws.Range("L75").ClearContents()

$ws.Range("H78").Value = 76839.164
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 76839.164
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = 230517.492
$ws.Range("N78").Value = -239253.492
$ws.Range("L78").ClearContents()

$ws.Range("H86").Value = 100314
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H88").Value = 1102.4166
$ws.Range("J88").Value = 1687.5
$ws.Range("L88").Value = 1687.5
$ws.Range("N88").Value = -2499.5

$ws.Range("H89").Value = 100314
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H91").Value = 1102.4166
$ws.Range("J91").Value = 1687.5
$ws.Range("L91").Value = 1687.5
$ws.Range("N91").Value = -4495.5

$ws.Range("H102").Value = 3370.2222
$ws.Range("I102").Value = 2976.1428
$ws.Range("K102").Value = 2976.1428
$ws.Range("M102").Value = -1354.1428

$ws.Range("H110").Value = 2432.2222
$ws.Range("I110").Value = 2578
$ws.Range("J110").Value = 2250
$ws.Range("K110").Value = 2578
$ws.Range("L110").Value = 2250
$ws.Range("M110").Value = -533
$ws.Range("N110").Value = -6340

$ws.Range("H122").Value = 2050.8572
$ws.Range("I122").Value = 1476.0834
$ws.Range("K122").Value = 4428.2502
$ws.Range("M122").Value = -1978.2502

$ws.Range("H132").Value = 5684.9443
$ws.Range("I132").Value = 4478.1665
$ws.Range("K132").Value = 13434.4995
$ws.Range("M132").Value = -10904.4995

$ws.Range("H136").Value = 3130.7778
$ws.Range("I136").Value = 2273.5386
$ws.Range("K136").Value = 6820.6158
$ws.Range("M136").Value = -4270.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1829.4231
$ws.Range("I86").Value = 1887.7693
$ws.Range("J86").Value = 1771.0769
$ws.Range("K86").Value = 1887.7693
$ws.Range("L86").Value = 1771.0769
$ws.Range("M86").Value = -764.7692999999999
$ws.Range("N86").Value = -4017.0769

$ws.Range("H89").Value = 1829.4231
$ws.Range("I89").Value = 1887.7693
$ws.Range("J89").Value = 1771.0769
$ws.Range("K89").Value = 9438.8465
$ws.Range("L89").Value = 8855.3845
$ws.Range("M89").Value = -3822.8465
$ws.Range("N89").Value = -20087.3845

$ws.Range("H99").Value = 8639.5
$ws.Range("I99").Value = 4506.3076
$ws.Range("J99").Value = 16315.429
$ws.Range("K99").Value = 4506.3076
$ws.Range("L99").Value = 16315.429
$ws.Range("M99").Value = -3008.3076
$ws.Range("N99").Value = -19311.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3931.36
$ws.Range("I31").Value = 1628.8
$ws.Range("J31").Value = 5466.4
$ws.Range("K31").Value = 1628.8
$ws.Range("L31").Value = 5466.4
$ws.Range("M31").Value = -1333.8
$ws.Range("N31").Value = -6056.4

$ws.Range("H34").Value = 3931.36
$ws.Range("I34").Value = 1628.8
$ws.Range("J34").Value = 5466.4
$ws.Range("K34").Value = 1628.8
$ws.Range("L34").Value = 5466.4
$ws.Range("M34").Value = -1426.8
$ws.Range("N34").Value = -5870.4

$ws.Range("H99").Value = 2499.6667
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2499.6667
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = 2499.6667
$ws.Range("N99").Value = -5495.6667
$ws.Range("L99").ClearContents()

$ws.Range("H122").Value = 1310.8096
$ws.Range("I122").Value = 1076.35
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 3229.05
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -779.0499999999997
$ws.Range("N122").Value = -22900

$ws.Range("H126").Value = 2499.6667
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2499.6667
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = 7499.000100000001
$ws.Range("N126").Value = -12439.0001
$ws.Range("L126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 18182040
$ws.Range("I122").Value = 207.83333
$ws.Range("J122").Value = 40000240
$ws.Range("K122").Value = 1870.49997
$ws.Range("L122").Value = 360002160
$ws.Range("M122").Value = 579.5000300000002
$ws.Range("N122").Value = -360007060

$ws.Range("H134").Value = 963.9
$ws.Range("I134").Value = 963.9
$ws.Range("K134").Value = 2891.7
$ws.Range("M134").Value = 2178.3

$ws.Range("H137").Value = 3518.4736
$ws.Range("I137").Value = 1853.625
$ws.Range("J137").Value = 4729.273
$ws.Range("K137").Value = 5560.875
$ws.Range("L137").Value = 14187.819
$ws.Range("M137").Value = -460.875
$ws.Range("N137").Value = -24387.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 751988
$ws.Range("I80").Value = 1001820.8
$ws.Range("K80").Value = 1001820.8
$ws.Range("M80").Value = -1000822.8

$ws.Range("H83").Value = 751988
$ws.Range("I83").Value = 1001820.8
$ws.Range("K83").Value = 5009104
$ws.Range("M83").Value = -5004112

$ws.Range("H102").Value = 4999
$ws.Range("I102").Value = 4999
$ws.Range("K102").Value = 4999
$ws.Range("M102").Value = -3377

$ws.Range("H126").Value = 2254.5
$ws.Range("I126").Value = 2254.5
$ws.Range("K126").Value = 6763.5
$ws.Range("M126").Value = -4293.5

$ws.Range("H132").Value = 17252454
$ws.Range("I132").Value = 25652434
$ws.Range("J132").Value = 10384.947
$ws.Range("K132").Value = 76957302
$ws.Range("L132").Value = 31154.841
$ws.Range("M132").Value = -76954772
$ws.Range("N132").Value = -36214.841

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2328.7778
$ws.Range("I82").Value = 1163.6666
$ws.Range("J82").Value = 2911.3333
$ws.Range("K82").Value = 1163.6666
$ws.Range("L82").Value = 2911.3333
$ws.Range("M82").Value = -802.6666
$ws.Range("N82").Value = -3633.3333

$ws.Range("H85").Value = 2328.7778
$ws.Range("I85").Value = 1163.6666
$ws.Range("J85").Value = 2911.3333
$ws.Range("K85").Value = 1163.6666
$ws.Range("L85").Value = 2911.3333
$ws.Range("M85").Value = 84.33339999999998
$ws.Range("N85").Value = -5407.3333

$ws.Range("H87").Value = 88151.39999999999
$ws.Range("J87").Value = 88151.39999999999
$ws.Range("L87").Value = 88151.39999999999
$ws.Range("N87").Value = -90397.39999999999

$ws.Range("H90").Value = 88151.39999999999
$ws.Range("J90").Value = 88151.39999999999
$ws.Range("L90").Value = 264454.2
$ws.Range("N90").Value = -275686.2

$ws.Range("H122").Value = 5324.25
$ws.Range("I122").Value = 4884.385
$ws.Range("K122").Value = 14653.155
$ws.Range("M122").Value = -12203.155

$ws.Range("H132").Value = 2388.6155
$ws.Range("I132").Value = 2388.6155
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7165.8465
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -4635.8465
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4653.8184
$ws.Range("I62").Value = 3699.6
$ws.Range("J62").Value = 5449
$ws.Range("K62").Value = 3699.6
$ws.Range("L62").Value = 5449
$ws.Range("M62").Value = -3075.6
$ws.Range("N62").Value = -6697

$ws.Range("H65").Value = 4653.8184
$ws.Range("I65").Value = 3699.6
$ws.Range("J65").Value = 5449
$ws.Range("K65").Value = 18498
$ws.Range("L65").Value = 27245
$ws.Range("M65").Value = -15378
$ws.Range("N65").Value = -33485

$ws.Range("H122").Value = 3611.7856
$ws.Range("I122").Value = 3431.7368
$ws.Range("J122").Value = 5322.25
$ws.Range("K122").Value = 10295.2104
$ws.Range("L122").Value = 15966.75
$ws.Range("M122").Value = -7845.2104
$ws.Range("N122").Value = -20866.75
